$d = $word.ActiveDocument

# The page footer block ("Ver no Jupiter Salvar em pdf Salvar em docx" and the
# "(c) 2020 ..." copyright line), along with the blank paragraph that separated
# them from the bibliography text, was removed from the end of the document.
#
# Locate the last paragraph of the bibliography section ("...tema de Engenharia
# Bioquimica") as an anchor, then remove the blank paragraph plus the two text
# paragraphs that immediately follow it.

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tema de Engenharia Bioqu*mica*") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not find anchor paragraph ending in 'tema de Engenharia Bioquimica'"
}

$blankPara = $anchor.Next()
$jupiterPara = $blankPara.Next()
$copyrightPara = $jupiterPara.Next()

if ($jupiterPara.Range.Text -notlike "*Ver no Jupiter*") {
    throw "Unexpected paragraph where 'Ver no Jupiter...' text was expected: $($jupiterPara.Range.Text)"
}
if ($copyrightPara.Range.Text -notlike "*Powered by Jekyll*") {
    throw "Unexpected paragraph where copyright text was expected: $($copyrightPara.Range.Text)"
}

$deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
$deleteRange.Delete()
